$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: append the newly-collected trend samples instead of overwriting
# the existing rows. New rows 26-39 carry fresh timestamps (UTC) and their
# corresponding values.
$newRows = @(
    @("20 Jan 2024 - 17:52 UTC", 0),
    @("20 Jan 2024 - 18:07 UTC", -46),
    @("20 Jan 2024 - 18:08 UTC", -46),
    @("20 Jan 2024 - 18:09 UTC", -44),
    @("20 Jan 2024 - 18:10 UTC", -44),
    @("20 Jan 2024 - 18:11 UTC", -43),
    @("20 Jan 2024 - 18:12 UTC", 0),
    @("20 Jan 2024 - 18:13 UTC", 0),
    @("20 Jan 2024 - 18:14 UTC", 0),
    @("20 Jan 2024 - 18:15 UTC", 0),
    @("20 Jan 2024 - 18:16 UTC", 0),
    @("20 Jan 2024 - 18:17 UTC", 0),
    @("20 Jan 2024 - 18:18 UTC", 0),
    @("20 Jan 2024 - 18:19 UTC", 0)
)

$startRow = 26
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
